$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 3
$ws.Range("B3").Value = "74883c9adf428f76ddf53da0f1c3c64d64b4c54a343267cd4a75820fdf78ae2f"
$ws.Range("F3").Value = 0.9080082135523614
$ws.Range("G3").Value = 0.08194768214086647
$ws.Range("H3").Value = 0.835765885301489
$ws.Range("I3").Value = 0.02127974188131047
$ws.Range("J3").Value = 0.8443396226415094
$ws.Range("K3").Value = 0.00002238482165143815
$ws.Range("N3").Value = 2587.831536531448

# Row 4
$ws.Range("N4").Value = 2315.666812181473

# Row 6
$ws.Range("B6").Value = "374c2b4f742ec5b1b2b74b3342bc7d42b939fa3fc6c8ba20bd51c62de31b1571"
$ws.Range("F6").Value = 0.9512557257937134
$ws.Range("G6").Value = 0.03056370849534771
$ws.Range("H6").Value = 0.8292790800530738
$ws.Range("I6").Value = 0.0472828486835165
$ws.Range("J6").Value = 0.8341686320754716
$ws.Range("K6").Value = 0.00007341973911814681
$ws.Range("N6").Value = 3583.84148144722

# Row 7
$ws.Range("N7").Value = 3643.592126607895

# Row 9
$ws.Range("N9").Value = 202.2005605697632

# Row 10
$ws.Range("N10").Value = 716.4432606697083

# Row 12
$ws.Range("N12").Value = 207.0547118186951

# Row 13
$ws.Range("N13").Value = 627.1305425167084

# Row 15
$ws.Range("B15").Value = "b2f359849744bda5e16b62f6188940f3178d00190f2ded924ed0e868160f3ccf"
$ws.Range("F15").Value = 0.8410565338276181
$ws.Range("G15").Value = 0.01767893469330031
$ws.Range("H15").Value = 0.6518918918918919
$ws.Range("I15").Value = 0.1022627855441957
$ws.Range("J15").Value = 0.661987041036717
$ws.Range("K15").Value = 0.005509831637616397
$ws.Range("M15").Value = 5
$ws.Range("N15").Value = 1674.26643705368

# Row 16
$ws.Range("N16").Value = 1199.010909080505

# Row 18
$ws.Range("B18").Value = "e7152604380e46f5aa215a9cea001ffb07d2b2bbacadd016300c8d0935b7de28"
$ws.Range("F18").Value = 0.8674698795180723
$ws.Range("G18").Value = 0.01492458224650935
$ws.Range("H18").Value = 0.6572972972972972
$ws.Range("I18").Value = 0.06360705195462565
$ws.Range("J18").Value = 0.6555075593952484
$ws.Range("K18").Value = 0.003201659451659455
$ws.Range("M18").Value = 512
$ws.Range("N18").Value = 2079.936646223068

# Row 19
$ws.Range("N19").Value = 1336.322180986404

# Row 21
$ws.Range("B21").Value = "954e8324ddf6a4a85178c177b0842e472c8514995767110da7c632dbd2dcffb9"
$ws.Range("F21").Value = 0.9397404057916686
$ws.Range("G21").Value = 0.005223628884854518
$ws.Range("H21").Value = 0.9024710255849552
$ws.Range("I21").Value = 0.08174070253590099
$ws.Range("J21").Value = 0.8891560996939222
$ws.Range("K21").Value = 0.007498373906090993
$ws.Range("M21").Value = 2
$ws.Range("N21").Value = 2879.724896669388

# Row 22
$ws.Range("N22").Value = 1771.492316246033
